$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.285.56'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '2.601.80'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D5').Value = '575.77'
$ws.Range('E5').Value = '  +3.37%  '
$ws.Range('D6').Value = '142.93'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = '2.604.29'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '3.059.01'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '24.33'
$ws.Range('E15').Value = '  +3.88%  '
$ws.Range('D16').Value = '60.278.44'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = '2.604.56'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '11.37'
$ws.Range('E19').Value = '  +6.73%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '346.26'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').Value = '6.88'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').Value = '63.04'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = '8.02'
$ws.Range('E28').Value = '  +5.36%  '
$ws.Range('D29').Value = '0.0₃0797'
$ws.Range('E29').Value = '  +1.73%  '
$ws.Range('E30').Value = '  +10.05%  '
$ws.Range('D31').Value = '6.37'
$ws.Range('E31').Value = '  +3.30%  '
$ws.Range('D33').Value = '166.27'
$ws.Range('E33').Value = '  +4.83%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.27'
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.30'
$ws.Range('E36').Value = '  +9.16%  '
$ws.Range('E37').Value = '  +7.41%  '
$ws.Range('E38').Value = '  +6.78%  '
$ws.Range('D39').Value = '38.07'
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').Value = '313.67'
$ws.Range('E40').Value = '  +6.74%  '
$ws.Range('E41').Value = '  +4.63%  '
$ws.Range('D42').Value = '0.838'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('E43').Value = '  -3.50%  '
$ws.Range('D44').Value = '0.0994'
$ws.Range('E44').Value = '  +1.38%  '
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('E46').Value = '  +2.52%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.605'
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '0.0551'
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').Value = '5.00'
$ws.Range('E49').Value = '  +4.60%  '
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').Value = '19.97'
$ws.Range('E51').Value = '  +5.03%  '
